# Add a new inventory line (TI part "U17" / TPS61253AYFFT) right above the
# "Total" row of the Sheet1 inventory table, pushing the Total row from
# row 76 down to row 77, and extend the print area / dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Push the existing "Total" row (currently row 76) down to row 77.
#    Copy-to-destination preserves the cell formatting (styles/borders/number
#    formats) of the Total row instead of creating brand new style entries.
$ws.Range("A76:F76").Copy($ws.Range("A77:F77"))
# Re-apply the SUM formula (Copy of a formula cell can lose the formula in
# this engine, so set it explicitly) - it still sums the same row range.
$ws.Range("D77").Formula = "=SUM(D6:D75)"
$ws.Rows.Item(77).RowHeight = 15.65

# 2) Turn the now-vacated row 76 into the new parts-list row. Copy the
#    formatting from row 75 (an ordinary data row) so the new row gets the
#    same styling (fonts/borders/number formats) as the rest of the table.
$ws.Range("A75:F75").Copy($ws.Range("A76:F76"))
$ws.Rows.Item(76).RowHeight = 13.8

# 3) Fill in the new part's data.
$ws.Range("A76").Value = "U17"
$ws.Range("B76").Value = 25
$ws.Range("C76").Value = 159.37
$ws.Range("D76").Value = 3984
$ws.Range("E76").Value = "TPS61253AYFFT"
$ws.Range("F76").Value = "TPS6125xA 3.8-MHz, 5-V / 4-A Boost in 1.2-mm x 1.3-mm WCSP"

# 4) Grow the print area so it still covers the whole table (now through
#    row 77 instead of row 76).
$ws.PageSetup.PrintArea = "`$A`$6:`$F`$77"

# 5) Best-effort: keep the view scrolled near the bottom of the table and
#    select the new last data cell, matching where the sheet was left.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 51
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F76").Select() | Out-Null
